$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing F (and some G) column values ---
$updates = @(
    @{ Row = 353; F = 725581 },
    @{ Row = 471; F = 67037 },
    @{ Row = 472; F = 51961 },
    @{ Row = 473; F = 40019; G = 41 },
    @{ Row = 476; F = 37376 },
    @{ Row = 477; F = 37209 },
    @{ Row = 480; F = 33566 },
    @{ Row = 484; F = 8277 },
    @{ Row = 485; F = 14015 },
    @{ Row = 486; F = 8949 },
    @{ Row = 487; F = 6811 },
    @{ Row = 497; F = 7785 },
    @{ Row = 498; F = 9199 },
    @{ Row = 499; F = 11472 },
    @{ Row = 500; F = 7785 },
    @{ Row = 502; F = 10613 },
    @{ Row = 503; F = 7496 },
    @{ Row = 504; F = 7524 },
    @{ Row = 505; F = 8562 },
    @{ Row = 506; F = 10865 },
    @{ Row = 507; F = 7245 },
    @{ Row = 508; F = 5705 },
    @{ Row = 509; F = 9663 },
    @{ Row = 510; F = 7900 },
    @{ Row = 511; F = 6841 },
    @{ Row = 512; F = 8530 },
    @{ Row = 513; F = 10421 },
    @{ Row = 514; F = 7012 },
    @{ Row = 515; F = 5044 },
    @{ Row = 516; F = 9398 },
    @{ Row = 517; F = 6785 },
    @{ Row = 518; F = 7127 },
    @{ Row = 519; F = 7914 },
    @{ Row = 520; F = 10196 },
    @{ Row = 521; F = 6718 },
    @{ Row = 522; F = 5022 },
    @{ Row = 523; F = 10086 },
    @{ Row = 524; F = 7759 },
    @{ Row = 525; F = 7479; G = 22 },
    @{ Row = 526; F = 8569; G = 26 },
    @{ Row = 527; F = 11085; G = 33 },
    @{ Row = 528; F = 7691; G = 20 },
    @{ Row = 529; F = 5216; G = 21 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    if ($u.ContainsKey('G')) {
        $ws.Cells.Item($u.Row, 7).Value = $u.G
    }
}

# --- Append new rows 530-532 ---
$newRows = @(
    @{ Row = 530; A = 44424; B = 393610; C = 7257; D = 74;  E = 12546; F = 11922; G = 39 },
    @{ Row = 531; A = 44425; B = 393722; C = 6472; D = 112; E = 12547; F = 8268;  G = 24 },
    @{ Row = 532; A = 44426; B = 393799; C = 6132; D = 77;  E = 12547; F = 6741;  G = 29 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
